$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vm_pu values per row (bus index), columns B,C,D,E,F,I,J,K,L,M,N
$data = @{
    2 = @{ "B" = 1.02; "C" = 1.016592233409688; "D" = 1.046399119885168; "E" = 1.018091725468495; "F" = 1.050215559784922; "I" = 1.039876406328814; "J" = 1.021811164796023; "K" = 1.049164665093498; "L" = 1.020938970791996; "M" = 1.052970455550394; "N" = 1.023262253853412 }
    3 = @{ "B" = 1.02; "C" = 1.017684158783479; "D" = 1.047013892584638; "E" = 1.019021475234721; "F" = 1.051045752225476; "I" = 1.040019172672328; "J" = 1.022537859145252; "K" = 1.049591584220517; "L" = 1.021674114338427; "M" = 1.053613003963875; "N" = 1.02398998019197 }
    4 = @{ "B" = 1.02; "C" = 1.018390804813626; "D" = 1.047409894452724; "E" = 1.019623551184558; "F" = 1.051581289682239; "I" = 1.040108996520666; "J" = 1.023007667958167; "K" = 1.049865343060104; "L" = 1.022149648130539; "M" = 1.054026479153831; "N" = 1.024460456187288 }
    5 = @{ "B" = 1.02; "C" = 1.018687902508343; "D" = 1.047575941350322; "E" = 1.019876775035381; "F" = 1.051806032084024; "I" = 1.040146145551195; "J" = 1.023205076988033; "K" = 1.049979834290862; "L" = 1.022349525603146; "M" = 1.054199752474686; "N" = 1.024658145560619 }
    6 = @{ "B" = 1.02; "C" = 1.018737787920443; "D" = 1.047603795945507; "E" = 1.019919298939722; "F" = 1.051843743953523; "I" = 1.040152347062297; "J" = 1.023238217045324; "K" = 1.049999022813523; "L" = 1.022383083737187; "M" = 1.054228813424589; "N" = 1.024691332680592 }
    7 = @{ "B" = 1.02; "C" = 1.01839477455393; "D" = 1.047412114880756; "E" = 1.019626934337894; "F" = 1.051584294265813; "I" = 1.040109495317633; "J" = 1.023010306134971; "K" = 1.049866875245561; "L" = 1.02215231904857; "M" = 1.054028796610451; "N" = 1.024463098110605 }
    8 = @{ "B" = 1.02; "C" = 1.016961234804202; "D" = 1.046607256098681; "E" = 1.018405842116494; "F" = 1.050496468239427; "I" = 1.039925183355141; "J" = 1.022056839543905; "K" = 1.049309458196892; "L" = 1.02118744757802; "M" = 1.053188082226355; "N" = 1.023508277487615 }
    9 = @{ "B" = 1.02; "C" = 1.014435880714035; "D" = 1.045175316195655; "E" = 1.016257703481427; "F" = 1.048567005899337; "I" = 1.039580878138876; "J" = 1.02037356948356; "K" = 1.04830825057775; "L" = 1.019486057370507; "M" = 1.051689126189499; "N" = 1.021822616990711 }
    10 = @{ "B" = 1.02; "C" = 1.01275276950033; "D" = 1.044211610873531; "E" = 1.014828037125288; "F" = 1.047272351621998; "I" = 1.039338277735567; "J" = 1.019249285977948; "K" = 1.047628124622301; "L" = 1.018351026230238; "M" = 1.050678153028175; "N" = 1.020696736873564 }
    11 = @{ "B" = 1.02; "C" = 1.012024064724213; "D" = 1.043792186110256; "E" = 1.0142095540969; "F" = 1.046709794097; "I" = 1.039230143347377; "J" = 1.018761960631019; "K" = 1.047330641086496; "L" = 1.01785936387489; "M" = 1.050237643508854; "N" = 1.020208719468745 }
    12 = @{ "B" = 1.02; "C" = 1.011753404588872; "D" = 1.043636074284094; "E" = 1.013979908022441; "F" = 1.046500541790176; "I" = 1.039189514618835; "J" = 1.018580870657714; "K" = 1.04721969557657; "L" = 1.017676710611242; "M" = 1.050073606841875; "N" = 1.020027372326906 }
    13 = @{ "B" = 1.02; "C" = 1.01181146151904; "D" = 1.043669575198045; "E" = 1.014029163988489; "F" = 1.046545440371238; "I" = 1.039198250546642; "J" = 1.018619718495444; "K" = 1.047243513963374; "L" = 1.017715891621802; "M" = 1.050108811874698; "N" = 1.020066275333023 }
    14 = @{ "B" = 1.02; "C" = 1.012001691606943; "D" = 1.043779288350935; "E" = 1.014190569724399; "F" = 1.046692503211909; "I" = 1.039226794399565; "J" = 1.018746993220683; "K" = 1.047321479403623; "L" = 1.017844266269651; "M" = 1.050224092578852; "N" = 1.02019373080297 }
    15 = @{ "B" = 1.02; "C" = 1.012118900376668; "D" = 1.043846844055812; "E" = 1.014290028533039; "F" = 1.046783074652233; "I" = 1.039244319911912; "J" = 1.018825401354808; "K" = 1.04736945731331; "L" = 1.017923358421051; "M" = 1.050295066301636; "N" = 1.020272250285639 }
    16 = @{ "B" = 1.02; "C" = 1.012801133139458; "D" = 1.044239401908507; "E" = 1.014869095868538; "F" = 1.047309645479314; "I" = 1.039345389305304; "J" = 1.019281617551302; "K" = 1.047647804874221; "L" = 1.018383652317206; "M" = 1.050707330357749; "N" = 1.02072911436146 }
    17 = @{ "B" = 1.02; "C" = 1.013229104206445; "D" = 1.044485073038111; "E" = 1.015232482729614; "F" = 1.047639424873533; "I" = 1.039407961720031; "J" = 1.0195676551756; "K" = 1.047821606983239; "L" = 1.018672332639854; "M" = 1.050965197173338; "N" = 1.021015558191994 }
    18 = @{ "B" = 1.02; "C" = 1.013478741561827; "D" = 1.04462816283169; "E" = 1.015444495312845; "F" = 1.047831590040868; "I" = 1.039444161255199; "J" = 1.019734447537323; "K" = 1.047922694751575; "L" = 1.018840697003303; "M" = 1.051115340944436; "N" = 1.021182587418002 }
    19 = @{ "B" = 1.02; "C" = 1.013563863011393; "D" = 1.044676917753278; "E" = 1.015516795479357; "F" = 1.047897081191961; "I" = 1.039456453803535; "J" = 1.019791311189561; "K" = 1.047957114161657; "L" = 1.018898101848859; "M" = 1.051166490995493; "N" = 1.021239531823149 }
    20 = @{ "B" = 1.02; "C" = 1.013183186028239; "D" = 1.044458736147232; "E" = 1.015193489050353; "F" = 1.047604062262083; "I" = 1.039401279097841; "J" = 1.019536971066987; "K" = 1.047802989446251; "L" = 1.018641361833183; "M" = 1.050937557963131; "N" = 1.020984830508428 }
    21 = @{ "B" = 1.02; "C" = 1.011945673227555; "D" = 1.043746989347804; "E" = 1.014143037381385; "F" = 1.046649204956991; "I" = 1.039218401714336; "J" = 1.01870951607481; "K" = 1.047298532848794; "L" = 1.017806463907526; "M" = 1.0501901566406; "N" = 1.02015620043525 }
    22 = @{ "B" = 1.02; "C" = 1.011167676165248; "D" = 1.043297642315653; "E" = 1.013483074704045; "F" = 1.046047151564378; "I" = 1.039100742273849; "J" = 1.018188824614912; "K" = 1.046978776183284; "L" = 1.017281368607567; "M" = 1.049717853862314; "N" = 1.019634769533756 }
    23 = @{ "B" = 1.02; "C" = 1.011580100152948; "D" = 1.043536023831144; "E" = 1.013832886160887; "F" = 1.046366471752866; "I" = 1.039163369220234; "J" = 1.018464894512314; "K" = 1.047148529753238; "L" = 1.017559746893781; "M" = 1.0499684557662; "N" = 1.019911231482078 }
    24 = @{ "B" = 1.02; "C" = 1.013203934448087; "D" = 1.044470637291208; "E" = 1.015211108444935; "F" = 1.047620041689781; "I" = 1.03940429960763; "J" = 1.019550836044723; "K" = 1.047811402799935; "L" = 1.018655356263027; "M" = 1.050950047753101; "N" = 1.020998715176022 }
    25 = @{ "B" = 1.02; "C" = 1.015088662302505; "D" = 1.045547115545593; "E" = 1.016812622475474; "F" = 1.049067296981843; "I" = 1.039672196386737; "J" = 1.020809106430874; "K" = 1.048569324225228; "L" = 1.019926044209377; "M" = 1.052078706468229; "N" = 1.022258772450447 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}